$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the new measurement types "Holden" and "Rizzie Spiral"
$ws.Rows("4:5").Insert()

# Copy formatting (bold, centered, bordered) from the row above onto the new rows
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New row 4: Holden
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.168608583874831
$ws.Range("D4").Value = 0.8886061077190477
$ws.Range("E4").Value = 0.9066785864273847
$ws.Range("F4").Value = 0.9066785864273847
$ws.Range("G4").Value = 0.8421302638119609
$ws.Range("H4").Value = 0.9154453387408147
$ws.Range("I4").Value = 1.006248355318265
$ws.Range("J4").Value = 0.7877817715820432
$ws.Range("K4").Value = 0.9066785864273847
$ws.Range("L4").Value = 0.7877817715820432
$ws.Range("M4").Value = 1.529053262247947
$ws.Range("N4").Value = 0.9066785864273847
$ws.Range("O4").Value = 1.529053262247947
$ws.Range("P4").Value = 1.158417516914995
$ws.Range("Q4").Value = 1.208829684983497
$ws.Range("R4").Value = 1.074504540085792
$ws.Range("S4").Value = 1.068480380516346
$ws.Range("T4").Value = 1.074504540085792
$ws.Range("U4").Value = 1.028029931994106
$ws.Range("V4").Value = 1.003759662880761
$ws.Range("W4").Value = 1.005569033715287

# New row 5: Rizzie Spiral
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 0.9749778847816253
$ws.Range("D5").Value = 0.982675016769275
$ws.Range("E5").Value = 2.839309824320281
$ws.Range("F5").Value = 2.839309824320281
$ws.Range("G5").Value = 1.131194740408092
$ws.Range("H5").Value = 1.216814230196834
$ws.Range("I5").Value = 0.2382745018627646
$ws.Range("J5").Value = 1.964527147187833
$ws.Range("K5").Value = 2.839309824320281
$ws.Range("L5").Value = 1.964527147187833
$ws.Range("M5").Value = 0.5957636841742221
$ws.Range("N5").Value = 2.839309824320281
$ws.Range("O5").Value = 0.5957636841742221
$ws.Range("P5").Value = 1.280145415681027
$ws.Range("Q5").Value = 0.7892193504717486
$ws.Range("R5").Value = 1.799866885227446
$ws.Range("S5").Value = 1.180988616043776
$ws.Range("T5").Value = 1.799866885227445
$ws.Range("U5").Value = 1.595568918112903
$ws.Range("V5").Value = 1.844317099354378
$ws.Range("W5").Value = 1.242942128712616

# Rename "Thomas Hex" -> "Matthies Hex" (now on row 11 after the insert shift)
$ws.Range("B11").Value = "Matthies Hex"

